$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks up front; row-insert does not re-anchor them to
# their shifted rows, so we rebuild the whole collection afterwards instead.
$ws.Hyperlinks.Delete()

# Insert a new row above the current row 2 (the most-recent-date data row),
# shifting all existing data rows (2-35) down to (3-36).
$ws.Rows("2:2").Insert()

# Copy formatting from the row below (the one that used to be row 2, now
# row 3) into the freshly inserted row 2 so it matches the table's styling.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the newly inserted row 2 with the latest price entry.
$ws.Cells.Item(2, 1).Value = 35
$ws.Cells.Item(2, 2).Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Cells.Item(2, 3).Value = "P1020"
$ws.Cells.Item(2, 4).Value = 267.75
$ws.Cells.Item(2, 5).Value = "23.09.2025"
$ws.Cells.Item(2, 6).Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-september-2025.pdf"

# Rebuild the hyperlinks for F2:F16 (rows with a Circular Link) in order.
$links = @(
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-17-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf",
    "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf"
)

for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Range("F$row"), $links[$i])
}

# Hyperlinks.Add reassigns the built-in "Hyperlink" style (blue/underline) to
# each touched cell; restore the plain centered style used throughout column F
# by pasting formats back in from the neighbouring (unstyled) column E.
$ws.Range("E2:E16").Copy()
$ws.Range("F2:F16").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
